# Updated cryptos list on Fri Jul 21 02:55:52 UTC 2023 with GitHub Actions
#
# Applies the latest price/volume(1h) refresh to the cryptos worksheet:
# column D ("Price") and column E ("Volume(1h)") are updated per-row for
# the coins whose quote moved since the last snapshot. Row 47 (PaxDollar)
# is unchanged this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "29.899.11"; E = "  -0.31%  " },
    @{ Row = 3; D = "1.898.64"; E = "  +0.05%  " },
    @{ Row = 4; D = $null; E = "  +0.06%  " },
    @{ Row = 5; D = $null; E = "  -5.27%  " },
    @{ Row = 6; D = "244.54"; E = "  +1.07%  " },
    @{ Row = 7; D = $null; E = "  +0.10%  " },
    @{ Row = 8; D = "0.3179"; E = "  -3.93%  " },
    @{ Row = 9; D = "25.55"; E = "  -4.42%  " },
    @{ Row = 10; D = "0.07181"; E = "  +1.57%  " },
    @{ Row = 11; D = "0.08116"; E = "  +0.37%  " },
    @{ Row = 12; D = "0.7724"; E = "  +1.74%  " },
    @{ Row = 13; D = "5.599"; E = "  +6.21%  " },
    @{ Row = 14; D = "1.882.91"; E = "  -0.81%  " },
    @{ Row = 15; D = "92.97"; E = "  +0.62%  " },
    @{ Row = 16; D = "6.185"; E = "  +4.99%  " },
    @{ Row = 17; D = "29.904.44"; E = "  -0.30%  " },
    @{ Row = 18; D = "13.98"; E = "  -1.05%  " },
    @{ Row = 19; D = "245.86"; E = "  +0.28%  " },
    @{ Row = 20; D = "0.000007772"; E = "  -0.20%  " },
    @{ Row = 21; D = "8.298"; E = "  +18.58%  " },
    @{ Row = 22; D = $null; E = "  +0.13%  " },
    @{ Row = 23; D = "2.148.88"; E = "  -0.03%  " },
    @{ Row = 24; D = $null; E = "  +0.16%  " },
    @{ Row = 25; D = "0.1676"; E = "  -4.66%  " },
    @{ Row = 26; D = "9.477"; E = "  +2.28%  " },
    @{ Row = 27; D = "164.34"; E = "  -1.13%  " },
    @{ Row = 28; D = $null; E = "  -0.87%  " },
    @{ Row = 29; D = "2.075"; E = "  -1.75%  " },
    @{ Row = 30; D = "1.410"; E = "  +3.58%  " },
    @{ Row = 31; D = "1.550"; E = "  +1.99%  " },
    @{ Row = 32; D = "4.511"; E = "  +4.87%  " },
    @{ Row = 33; D = "0.05639"; E = "  -3.30%  " },
    @{ Row = 34; D = "4.088"; E = "  +0.14%  " },
    @{ Row = 35; D = $null; E = "  +1.14%  " },
    @{ Row = 36; D = "0.7451"; E = "  +1.61%  " },
    @{ Row = 37; D = "1.003"; E = "  +0.47%  " },
    @{ Row = 38; D = "2.637"; E = "  -3.15%  " },
    @{ Row = 39; D = "0.01940"; E = "  +1.01%  " },
    @{ Row = 40; D = "2.789"; E = "  +0.56%  " },
    @{ Row = 41; D = "1.174.93"; E = "  +16.27%  " },
    @{ Row = 42; D = "74.88"; E = "  +3.13%  " },
    @{ Row = 43; D = "0.4441"; E = "  -0.25%  " },
    @{ Row = 44; D = "5.960"; E = "  +1.31%  " },
    @{ Row = 45; D = "0.8560"; E = "  +1.55%  " },
    @{ Row = 46; D = "104.69"; E = "  +2.78%  " },
    @{ Row = 48; D = "10.14"; E = "  +3.37%  " },
    @{ Row = 49; D = "1.888"; E = "  -0.09%  " },
    @{ Row = 50; D = "7.503"; E = "  -1.18%  " },
    @{ Row = 51; D = "2.980"; E = "  +9.62%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $origStyle = $cellD.Style
        # Force text storage so values like "1.410" / "0.000007772" keep their
        # literal digits/trailing zeros instead of being coerced to numbers,
        # then restore the original (default) style so no formatting drifts.
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = $origStyle
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
